# Apply updated dni_ciu (column E) and recompute PORC_AVANCE (column G)
# values for the affected rows in the "tabla_desagregada_mcp_merged" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new value for column E (dni_ciu)
$updates = @{
    2  = 1147
    3  = 870
    4  = 1040
    6  = 1039
    7  = 1091
    8  = 863
    9  = 1351
    10 = 1009
    11 = 1704
    12 = 793
    14 = 2136
}

foreach ($row in $updates.Keys) {
    $newE = $updates[$row]
    $ws.Cells.Item($row, 5).Value = $newE

    $pop = $ws.Cells.Item($row, 4).Value2
    $ws.Cells.Item($row, 7).Value = ($newE / $pop) * 100
}
